$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the existing hyperlinks (A1 and A2:A3) - cells keep their text,
# but are no longer clickable links, and lose the "Hyperlink" cell style.
$ws.Cells.Hyperlinks.Delete()

# Update / add the URL text values.
$ws.Range("A1").Value = "https://www.microsoft.com/en-us/surface"
$ws.Range("A2").Value = "https://www.microsoft.com/en-us/surface/devices/surface-pro"
$ws.Range("A3").Value = "https://www.microsoft.com/en-us/surface/devices/surface-laptop"
$ws.Range("A4").Value = "microsoft.com/en-us/surface/devices/compare-devices"

# Strip the now-unused Hyperlink formatting back to the default Normal style.
$ws.Range("A1:A4").Style = "Normal"

# Resize column A and move the active selection, matching the refreshed view.
$ws.Columns("A").ColumnWidth = 64.7
$ws.Range("B9").Select() | Out-Null
